# Atualizacao de bases das ligas, do dia: 13-06-2024 as 19:35
# Swap the data (columns B:AD) between three pairs of rows that were
# re-ordered upstream; column A (the sequential id) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 185 <-> row 186 (columns B:AD), keep column A (id) unchanged
$ws.Range("B185").Value = 6810162
$ws.Range("C185").Value = "Belgium First Division A"
$ws.Range("D185").Value = 45322.70833333334
$ws.Range("E185").Value = "Standard Liege"
$ws.Range("F185").Value = "Antwerp"
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1
$ws.Range("I185").Value = 0
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = "A"
$ws.Range("L185").Value = 4
$ws.Range("M185").Value = 3.6
$ws.Range("N185").Value = 1.85
$ws.Range("O185").Value = 3.1
$ws.Range("P185").Value = 3.2
$ws.Range("Q185").Value = 2.3
$ws.Range("R185").Value = 0.25
$ws.Range("S185").Value = 1.8
$ws.Range("T185").Value = 2.05
$ws.Range("U185").Value = 2.25
$ws.Range("V185").Value = 1.875
$ws.Range("W185").Value = 1.975
$ws.Range("X185").Value = -1
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = 1.3
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = 1.05
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.9750000000000001
$ws.Range("B186").Value = 6810164
$ws.Range("C186").Value = "Belgium First Division A"
$ws.Range("D186").Value = 45322.70833333334
$ws.Range("E186").Value = "Union Saint Gilloise"
$ws.Range("F186").Value = "RWD Molenbeek"
$ws.Range("G186").Value = 3
$ws.Range("H186").Value = 2
$ws.Range("I186").Value = 2
$ws.Range("J186").Value = 0
$ws.Range("K186").Value = "H"
$ws.Range("L186").Value = 1.2
$ws.Range("M186").Value = 7
$ws.Range("N186").Value = 12
$ws.Range("O186").Value = 1.142
$ws.Range("P186").Value = 8.5
$ws.Range("Q186").Value = 15
$ws.Range("R186").Value = -2.25
$ws.Range("S186").Value = 1.925
$ws.Range("T186").Value = 1.925
$ws.Range("U186").Value = 3.5
$ws.Range("V186").Value = 2.025
$ws.Range("W186").Value = 1.825
$ws.Range("X186").Value = 0.1419999999999999
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 0.925
$ws.Range("AC186").Value = 1.025
$ws.Range("AD186").Value = -1

# Swap row 187 <-> row 188 (columns B:AD), keep column A (id) unchanged
$ws.Range("B187").Value = 6810163
$ws.Range("C187").Value = "Belgium First Division A"
$ws.Range("D187").Value = 45323.6875
$ws.Range("E187").Value = "SintTruidense"
$ws.Range("F187").Value = "Gent"
$ws.Range("G187").Value = 4
$ws.Range("H187").Value = 1
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = 1
$ws.Range("K187").Value = "H"
$ws.Range("L187").Value = 3.6
$ws.Range("M187").Value = 3.6
$ws.Range("N187").Value = 1.95
$ws.Range("O187").Value = 3.25
$ws.Range("P187").Value = 3.4
$ws.Range("Q187").Value = 2.15
$ws.Range("R187").Value = 0.25
$ws.Range("S187").Value = 1.95
$ws.Range("T187").Value = 1.9
$ws.Range("U187").Value = 2.5
$ws.Range("V187").Value = 1.975
$ws.Range("W187").Value = 1.875
$ws.Range("X187").Value = 2.25
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 0.95
$ws.Range("AB187").Value = -1
$ws.Range("AC187").Value = 0.9750000000000001
$ws.Range("AD187").Value = -1
$ws.Range("B188").Value = 6810166
$ws.Range("C188").Value = "Belgium First Division A"
$ws.Range("D188").Value = 45323.6875
$ws.Range("E188").Value = "KV Mechelen"
$ws.Range("F188").Value = "Anderlecht"
$ws.Range("G188").Value = 2
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = 1
$ws.Range("K188").Value = "D"
$ws.Range("L188").Value = 3.5
$ws.Range("M188").Value = 3.5
$ws.Range("N188").Value = 2
$ws.Range("O188").Value = 3
$ws.Range("P188").Value = 3.5
$ws.Range("Q188").Value = 2.2
$ws.Range("R188").Value = 0.25
$ws.Range("S188").Value = 1.925
$ws.Range("T188").Value = 1.925
$ws.Range("U188").Value = 2.5
$ws.Range("V188").Value = 1.875
$ws.Range("W188").Value = 1.975
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = 2.5
$ws.Range("Z188").Value = -1
$ws.Range("AA188").Value = 0.4625
$ws.Range("AB188").Value = -0.5
$ws.Range("AC188").Value = 0.875
$ws.Range("AD188").Value = -1

# Swap row 278 <-> row 279 (columns B:AD), keep column A (id) unchanged
$ws.Range("B278").Value = 7979473
$ws.Range("C278").Value = "Belgium First Division A"
$ws.Range("D278").Value = 45406.64583333334
$ws.Range("E278").Value = "Anderlecht"
$ws.Range("F278").Value = "Cercle Brugge"
$ws.Range("G278").Value = 3
$ws.Range("H278").Value = 0
$ws.Range("I278").Value = 2
$ws.Range("J278").Value = 0
$ws.Range("K278").Value = "H"
$ws.Range("L278").Value = 1.909
$ws.Range("M278").Value = 3.6
$ws.Range("N278").Value = 3.8
$ws.Range("O278").Value = 1.8
$ws.Range("P278").Value = 3.8
$ws.Range("Q278").Value = 4
$ws.Range("R278").Value = -0.5
$ws.Range("S278").Value = 1.85
$ws.Range("T278").Value = 2
$ws.Range("U278").Value = 2.75
$ws.Range("V278").Value = 1.85
$ws.Range("W278").Value = 2
$ws.Range("X278").Value = 0.8
$ws.Range("Y278").Value = -1
$ws.Range("Z278").Value = -1
$ws.Range("AA278").Value = 0.8500000000000001
$ws.Range("AB278").Value = -1
$ws.Range("AC278").Value = 0.425
$ws.Range("AD278").Value = -0.5
$ws.Range("B279").Value = 7979357
$ws.Range("C279").Value = "Belgium First Division A"
$ws.Range("D279").Value = 45406.64583333334
$ws.Range("E279").Value = "Club Brugge"
$ws.Range("F279").Value = "Genk"
$ws.Range("G279").Value = 4
$ws.Range("H279").Value = 0
$ws.Range("I279").Value = 1
$ws.Range("J279").Value = 0
$ws.Range("K279").Value = "H"
$ws.Range("L279").Value = 1.85
$ws.Range("M279").Value = 3.75
$ws.Range("N279").Value = 3.9
$ws.Range("O279").Value = 1.75
$ws.Range("P279").Value = 3.75
$ws.Range("Q279").Value = 4.5
$ws.Range("R279").Value = -0.75
$ws.Range("S279").Value = 2
$ws.Range("T279").Value = 1.85
$ws.Range("U279").Value = 2.75
$ws.Range("V279").Value = 2.025
$ws.Range("W279").Value = 1.825
$ws.Range("X279").Value = 0.75
$ws.Range("Y279").Value = -1
$ws.Range("Z279").Value = -1
$ws.Range("AA279").Value = 1
$ws.Range("AB279").Value = -1
$ws.Range("AC279").Value = 1.025
$ws.Range("AD279").Value = -1
